$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the last used row in column A (the "Gold data" table starts at row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# New data: date label and the price description text
$ws.Cells.Item($newRow, 1).Value = "24-12-2025"
$ws.Cells.Item($newRow, 2).Value = "The price of gold in India today is ₹13,893 per gram for 24 karat gold, ₹12,735 per gram for 22 karat gold and ₹10,420 per gram for 18 karat gold (also called 999 gold)."

# Match the style (border + alignment) of the previous row
$ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow, 2)).Copy()
$ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 2)).PasteSpecial(-4122)
$excel.CutCopyMode = 0
